$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmt = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# ---------------------------------------------------------------
# Row 1 - headers: rename D1 and add new header cells E1:G1
# (write in this order so new shared-strings land in the same
#  sequence as the target workbook)
# ---------------------------------------------------------------
$ws.Range("D1").Value = "Quantity on Board"
$ws.Range("E1").Value = "Order"
$ws.Range("F1").Value = "Unit Cost"

# ---------------------------------------------------------------
# Row 2 - T1,2 / N MOSFET
# ---------------------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0.3
$ws.Range("G2").Formula = "=E2*F2"
$ws.Range("H2").Value = "Verified pinout"

# ---------------------------------------------------------------
# Row 3 - D4 / Power diode
# ---------------------------------------------------------------
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.38
$ws.Range("G3:G12").Formula = "=E3*F3"

# ---------------------------------------------------------------
# Row 4 - F1 / 1.5 A fuse
# ---------------------------------------------------------------
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.21
$ws.Range("H4").Value = "From LIDAR board"

# ---------------------------------------------------------------
# Row 5 - R / 1K 0603 Resistor
# ---------------------------------------------------------------
$ws.Range("C5").Value = "541-3991-1-ND"
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 0.022
$ws.Range("H5").Value = "From LIDAR board"
$ws.Range("I5").Value = "Buying 10 is cheaper than buying 3."

# ---------------------------------------------------------------
# Row 6 - 10K
# ---------------------------------------------------------------
$ws.Range("E6").Value = 0
$ws.Range("F6").NumberFormat = $fmt

# ---------------------------------------------------------------
# Row 7 - 5V REG
# ---------------------------------------------------------------
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1.88
$ws.Range("H7").Value = "From LIDAR board"

# ---------------------------------------------------------------
# Row 8 - 10 uF Cap
# ---------------------------------------------------------------
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.27

# ---------------------------------------------------------------
# Row 9 - 47 uF Ceramic Cap
# ---------------------------------------------------------------
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6

# ---------------------------------------------------------------
# Row 10 - Reset Button
# ---------------------------------------------------------------
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.52

# ---------------------------------------------------------------
# Row 11 - Signal Diodes
# ---------------------------------------------------------------
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 0.14
$ws.Range("H11").Value = "Not SBAS16HT1GOSCT-ND"

# ---------------------------------------------------------------
# Row 12 - Green SMD LED
# ---------------------------------------------------------------
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.29

# ---------------------------------------------------------------
# Row 1 - last new header (after the rest so "Total Cost" lands
# at the end of the shared-strings table, matching the target)
# ---------------------------------------------------------------
$ws.Range("G1").Value = "Total Cost"

# ---------------------------------------------------------------
# Row 14 - grand total
# ---------------------------------------------------------------
$ws.Range("G14").Formula = "=SUM(G2:G12)"

# ---------------------------------------------------------------
# Currency formatting for the new Unit Cost / Total Cost columns
# ---------------------------------------------------------------
$ws.Range("F2:G12").NumberFormat = $fmt
$ws.Range("G14").NumberFormat = $fmt

# ---------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 15.666666666666666
$ws.Columns.Item(5).ColumnWidth = 6.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.333333333333333
$ws.Columns.Item(7).ColumnWidth = 8.333333333333334
$ws.Columns.Item(8).ColumnWidth = 15

# ---------------------------------------------------------------
# Selection
# ---------------------------------------------------------------
$ws.Range("G15").Select() | Out-Null
